$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.909.35'
$ws.Range('D3').Value = '2.399.36'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '560.89'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '142.29'
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = '5.25'
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '25.51'
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('E14').Value = '  -1.73%  '
$ws.Range('D15').Value = '2.833.80'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').Value = '61.832.25'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '2.399.52'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = '11.19'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').Value = '320.20'
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').Value = '65.81'
$ws.Range('D24').Value = '1.73'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').Value = '8.76'
$ws.Range('E25').Value = '  -4.29%  '
$ws.Range('D26').Value = '560.48'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').Value = '2.520.25'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').Value = '0.0₃0928'
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('E31').Value = '  -4.67%  '
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -4.92%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').Value = '4.70'
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').Value = '152.33'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('E38').Value = '  -5.85%  '
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('D40').Value = '18.53'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('E41').Value = '  -5.88%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '2.24'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('D44').Value = '147.07'
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('D45').Value = '3.59'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').Value = '0.0528'
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('D47').Value = '19.69'
$ws.Range('E47').Value = '  -3.80%  '
$ws.Range('D48').Value = '0.584'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('D49').Value = '0.0916'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('E51').Value = '  +0.39%  '
